$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells keep their exact text representation (avoid numeric auto-conversion)
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).NumberFormat = "@"

# Apply updated cell values
$ws.Cells.Item(2, 4).Value = "60.870.08"
$ws.Cells.Item(2, 5).Value = "  +3.08%  "
$ws.Cells.Item(3, 4).Value = "2.686.69"
$ws.Cells.Item(3, 5).Value = "  +0.99%  "
$ws.Cells.Item(4, 4).Value = "0.999"
$ws.Cells.Item(4, 5).Value = "  -0.23%  "
$ws.Cells.Item(5, 4).Value = "523.97"
$ws.Cells.Item(5, 5).Value = "  +1.34%  "
$ws.Cells.Item(6, 4).Value = "144.90"
$ws.Cells.Item(6, 5).Value = "  +0.28%  "
$ws.Cells.Item(7, 5).Value = "  +0.01%  "
$ws.Cells.Item(8, 5).Value = "  +1.63%  "
$ws.Cells.Item(9, 4).Value = "2.705.02"
$ws.Cells.Item(9, 5).Value = "  +0.69%  "
$ws.Cells.Item(10, 4).Value = "6.47"
$ws.Cells.Item(10, 5).Value = "  +3.42%  "
$ws.Cells.Item(11, 5).Value = "  +0.11%  "
$ws.Cells.Item(12, 5).Value = "  +0.33%  "
$ws.Cells.Item(13, 5).Value = "  +2.46%  "
$ws.Cells.Item(14, 4).Value = "3.159.16"
$ws.Cells.Item(14, 5).Value = "  +0.45%  "
$ws.Cells.Item(15, 4).Value = "60.771.71"
$ws.Cells.Item(15, 5).Value = "  +2.97%  "
$ws.Cells.Item(16, 4).Value = "21.32"
$ws.Cells.Item(16, 5).Value = "  +1.23%  "
$ws.Cells.Item(17, 5).Value = "  +0.59%  "
$ws.Cells.Item(18, 4).Value = "2.703.08"
$ws.Cells.Item(18, 5).Value = "  +0.72%  "
$ws.Cells.Item(19, 4).Value = "349.73"
$ws.Cells.Item(19, 5).Value = "  -1.09%  "
$ws.Cells.Item(20, 4).Value = "4.51"
$ws.Cells.Item(20, 5).Value = "  -0.75%  "
$ws.Cells.Item(21, 4).Value = "10.56"
$ws.Cells.Item(21, 5).Value = "  +1.33%  "
$ws.Cells.Item(22, 5).Value = "  +1.43%  "
$ws.Cells.Item(23, 4).Value = "0.997"
$ws.Cells.Item(23, 5).Value = "  +0.06%  "
$ws.Cells.Item(24, 4).Value = "63.86"
$ws.Cells.Item(24, 5).Value = "  +2.83%  "
$ws.Cells.Item(25, 4).Value = "0.421"
$ws.Cells.Item(25, 5).Value = "  +0.11%  "
$ws.Cells.Item(26, 5).Value = "  +5.04%  "
$ws.Cells.Item(27, 5).Value = "  +0.77%  "
$ws.Cells.Item(28, 4).Value = "0.0₃0819"
$ws.Cells.Item(28, 5).Value = "  +1.23%  "
$ws.Cells.Item(29, 4).Value = "7.33"
$ws.Cells.Item(29, 5).Value = "  +1.54%  "
$ws.Cells.Item(30, 4).Value = "6.87"
$ws.Cells.Item(30, 5).Value = "  +8.08%  "
$ws.Cells.Item(31, 5).Value = "  +0.08%  "
$ws.Cells.Item(32, 4).Value = "19.28"
$ws.Cells.Item(32, 5).Value = "  +1.04%  "
$ws.Cells.Item(33, 5).Value = "  +0.93%  "
$ws.Cells.Item(34, 4).Value = "150.29"
$ws.Cells.Item(34, 5).Value = "  -0.23%  "
$ws.Cells.Item(35, 5).Value = "  +5.62%  "
$ws.Cells.Item(36, 4).Value = "1.26"
$ws.Cells.Item(36, 5).Value = "  +9.88%  "
$ws.Cells.Item(37, 4).Value = "0.949"
$ws.Cells.Item(37, 5).Value = "  -2.24%  "
$ws.Cells.Item(38, 4).Value = "0.878"
$ws.Cells.Item(38, 5).Value = "  +3.96%  "
$ws.Cells.Item(39, 4).Value = "1.53"
$ws.Cells.Item(39, 5).Value = "  +7.88%  "
$ws.Cells.Item(40, 4).Value = "36.94"
$ws.Cells.Item(40, 5).Value = "  +0.76%  "
$ws.Cells.Item(41, 5).Value = "  -2.25%  "
$ws.Cells.Item(42, 4).Value = "282.49"
$ws.Cells.Item(42, 5).Value = "  +0.48%  "
$ws.Cells.Item(43, 4).Value = "20.11"
$ws.Cells.Item(43, 5).Value = "  +2.42%  "
$ws.Cells.Item(44, 2).Value = "Stellar"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(44, 4).Value = "0.0992"
$ws.Cells.Item(44, 5).Value = "  +0.24%  "
$ws.Cells.Item(45, 2).Value = "FirstDigitalUSD"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Cells.Item(45, 4).Value = "0.997"
$ws.Cells.Item(45, 5).Value = "  -0.01%  "
$ws.Cells.Item(46, 4).Value = "0.610"
$ws.Cells.Item(46, 5).Value = "  -1.50%  "
$ws.Cells.Item(47, 4).Value = "2.142.60"
$ws.Cells.Item(47, 5).Value = "  +6.69%  "
$ws.Cells.Item(48, 2).Value = "RenderToken"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(48, 4).Value = "4.92"
$ws.Cells.Item(48, 5).Value = "  +6.57%  "
$ws.Cells.Item(49, 2).Value = "Hedera"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(49, 4).Value = "0.0540"
$ws.Cells.Item(49, 5).Value = "  +1.52%  "
$ws.Cells.Item(50, 5).Value = "  +1.49%  "
$ws.Cells.Item(51, 5).Value = "  +1.60%  "
